$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Row 4: date value change (45125 -> 45156) and E4 ("" -> "/") ---
$ws1.Cells.Item(4, 1).Value = 45156
$ws1.Cells.Item(4, 5).Value = '/'

# --- Append new rows 21-32 to sheet "订单" ---
# Row 21
$ws1.Cells.Item(21, 1).Value = 45264
$ws1.Cells.Item(21, 2).Value = 20
$ws1.Cells.Item(21, 4).Value = '2021112501-原始数据提供'
$ws1.Cells.Item(21, 6).Value = '黄礼闯'
$ws1.Cells.Item(21, 7).Value = 45260
$ws1.Cells.Item(21, 8).Value = '完成'
$ws1.Cells.Item(21, 9).Value = '原始数据提供'

# Row 22
$ws1.Cells.Item(22, 1).Value = 45265
$ws1.Cells.Item(22, 2).Value = 21
$ws1.Cells.Item(22, 4).Value = 'Treg 细胞差异表达基因'
$ws1.Cells.Item(22, 6).Value = '黄礼闯'
$ws1.Cells.Item(22, 7).Value = 45264
$ws1.Cells.Item(22, 8).Value = '完成'
$ws1.Cells.Item(22, 9).Value = 'ccRCC 单细胞数据的 Treg 细胞差异表达基因'

# Row 23
$ws1.Cells.Item(23, 1).Value = 45272
$ws1.Cells.Item(23, 2).Value = 22
$ws1.Cells.Item(23, 3).Value = ' 01-订单编号：IN2023120404 02-区域-销售：付梓欣 03-上级主管：王立家 04-医院：浙江省人民医院 05-科室/职称：康复科 06-电话： 07-项目（确定A/B套餐）：生信分析 08-分值：sci 1.5-2(生信) 09-定题题目：无 10-时间要求：2025年5月31前 11-总价： 12-定金：已付 13-评估人员 ：陶安琪，孙慧 14-技术支持（沟通情况）：薛富才（1次）。客户做生信，与客户课题相关。15-附件：评估邮件汇总，临床实验方案 16-备注：16-1.客户分级（需要文章晋升，或者单纯课题结题。老客户。） 16-2.谈单承诺,（与客户课题具有相关性；不投中科院预警杂志。） 17-项目负责人：杨弘 客户：章玮 '
$ws1.Cells.Item(23, 4).Value = 'IN2023120404'
$ws1.Cells.Item(23, 5).Value = '1.5-2'
$ws1.Cells.Item(23, 6).Value = '黄礼闯'
$ws1.Cells.Item(23, 7).Value = 45264
$ws1.Cells.Item(23, 8).Value = '完成'
$ws1.Cells.Item(23, 9).Value = 'RNA-seq 探究 rTMS 对 SCI 和 NP 的影响'

# Row 24
$ws1.Cells.Item(24, 1).Value = 45266
$ws1.Cells.Item(24, 2).Value = 23
$ws1.Cells.Item(24, 3).Value = ' 01-订单编号：实验：S2023120402；毕业论文：N2023120403 02-区域-销售：吴晓凤、叶立欢 03-上级主管： 04-医院：富阳第一人民医院 05-科室/职称： 06-电话： 07-项目（确定A/B套餐）：实验+毕业论文 08-分值： 09-定题题目： 10-时间要求：2024年3月 11-总价： 12-定金：已付 13-评估人员 ：陈颖+吴晓凤 14-技术支持（沟通情况）：吴晓凤 15-附件：定金，报价，技术支持与客户沟通总结 16-备注：1.客户实验+毕业论文） 2.谈单承诺,（复方细胞实验只提供三次有效数据，动物实验的材料由我司代买，客户自己付钱，动物检测部分结束后会给蜡块及切片） 17-项目负责人：杨弘 '
$ws1.Cells.Item(24, 4).Value = '实验：S2023120402；毕业论文：N2023120403'
$ws1.Cells.Item(24, 6).Value = '黄礼闯'
$ws1.Cells.Item(24, 7).Value = 45266
$ws1.Cells.Item(24, 8).Value = '完成'
$ws1.Cells.Item(24, 9).Value = '补肾健脾汤网络药理学分析'

# Row 25
$ws1.Cells.Item(25, 1).Value = 45266
$ws1.Cells.Item(25, 2).Value = 24
$ws1.Cells.Item(25, 4).Value = '半夏泻心汤网络药理学分析'
$ws1.Cells.Item(25, 6).Value = '黄礼闯'
$ws1.Cells.Item(25, 7).Value = 45266
$ws1.Cells.Item(25, 8).Value = '完成'
$ws1.Cells.Item(25, 9).Value = '半夏泻心汤网络药理学分析'

# Row 26
$ws1.Cells.Item(26, 1).Value = 45268
$ws1.Cells.Item(26, 2).Value = 25
$ws1.Cells.Item(26, 4).Value = '方和敬-白茅根-IgA网络药理学分析'
$ws1.Cells.Item(26, 6).Value = '黄礼闯'
$ws1.Cells.Item(26, 7).Value = 45267
$ws1.Cells.Item(26, 8).Value = '完成'
$ws1.Cells.Item(26, 9).Value = '白茅根-IgA网络药理学分析'

# Row 27
$ws1.Cells.Item(27, 1).Value = 45275
$ws1.Cells.Item(27, 2).Value = 26
$ws1.Cells.Item(27, 4).Value = 'S2023110704'
$ws1.Cells.Item(27, 6).Value = '黄礼闯'
$ws1.Cells.Item(27, 7).Value = 45272
$ws1.Cells.Item(27, 8).Value = '完成'
$ws1.Cells.Item(27, 9).Value = 'RNA的结合位点'

# Row 28
$ws1.Cells.Item(28, 1).Value = 45273
$ws1.Cells.Item(28, 2).Value = 27
$ws1.Cells.Item(28, 4).Value = '陈云杰测序结果差异分析'
$ws1.Cells.Item(28, 6).Value = '黄礼闯'
$ws1.Cells.Item(28, 7).Value = 45272
$ws1.Cells.Item(28, 8).Value = '完成'
$ws1.Cells.Item(28, 9).Value = '测序结果差异分析'

# Row 29
$ws1.Cells.Item(29, 1).Value = 45274
$ws1.Cells.Item(29, 2).Value = 28
$ws1.Cells.Item(29, 4).Value = 'SN2023011001'
$ws1.Cells.Item(29, 6).Value = '黄礼闯'
$ws1.Cells.Item(29, 7).Value = 45274
$ws1.Cells.Item(29, 8).Value = '完成'
$ws1.Cells.Item(29, 9).Value = '陈云杰测序数据分析++'

# Row 30
$ws1.Cells.Item(30, 1).Value = 45274
$ws1.Cells.Item(30, 2).Value = 29
$ws1.Cells.Item(30, 3).Value = ' 01-订单编号： 02-区域-销售：江苏-郭树仁 03-上级主管：王立家 04-医院：浙江省中 05-科室/职称：耳鼻喉科 06-电话： 07-项目（确定A/B套餐）：A 08-分值：3-5分 中科院三区 09-定题题目： 10-时间要求：2023/11/23-2025/5/23（18月）1月15日给国青标书，24年6月份给中管局、卫生厅标书， 11-总价： 12-定金： 13-评估人员 ：吴晨 14-技术支持（沟通情况）：薛富才（4次），吴晨 15-附件：定金截图（包含标书），方案，合同，实验报价 16-备注：1.潜力客户，做项目为了后续拿课题。 2.谈单承诺,（沟通过程中有答应客户的要求请尽数附上）：1月中给国青标书，6月初给中管局、卫生厅标书，先做预实验，争取国青本子中包含一点数据结果（和吴晓凤经理沟通过） 3、实验分阶段进行，除了与课题相关的实验部分，后续文章部分的实验等通知再进行。此外预实验分成两大模块，具体情况可与售前技术吴晨或者薛富才沟通。 4、所有实验需要走实验项目，需要实验分阶段汇报，动物实验保存图片和视频。 '
$ws1.Cells.Item(30, 4).Value = 'A2023112405'
$ws1.Cells.Item(30, 5).Value = 'sci3-5分 中科院三区'
$ws1.Cells.Item(30, 6).Value = '黄礼闯'
$ws1.Cells.Item(30, 7).Value = 45274
$ws1.Cells.Item(30, 8).Value = '完成'
$ws1.Cells.Item(30, 9).Value = '靳阳子生信支持业务'

# Row 31
$ws1.Cells.Item(31, 1).Value = 45285
$ws1.Cells.Item(31, 2).Value = 30
$ws1.Cells.Item(31, 3).Value = ' 01-订单编号：IN2023122103 02-区域-销售：周燕青 03-上级主管：柳叶 04-医院： 05-科室/职称：消化内科 06-电话： 07-项目（确定A/B套餐）：生信分析 08-分值： 09-定题题目： 10-时间要求：2023年12月27日前完成 11-总价： 12-定金：已付（结清） 13-评估人员 ：林婧羽  14-技术支持：林婧羽、薛富才（沟通情况）： 15-附件：后续合同再补充  '
$ws1.Cells.Item(31, 4).Value = 'IN2023122103'
$ws1.Cells.Item(31, 5).Value = '/'
$ws1.Cells.Item(31, 6).Value = '黄礼闯'
$ws1.Cells.Item(31, 7).Value = 45281
$ws1.Cells.Item(31, 8).Value = '完成'
$ws1.Cells.Item(31, 9).Value = '胆结石RNA-seq结合肠道菌、代谢物筛选关键差异表达基因'

# Row 32
$ws1.Cells.Item(32, 1).Value = 45285
$ws1.Cells.Item(32, 2).Value = 31
$ws1.Cells.Item(32, 4).Value = '周芳药方-草药-单体-靶点'
$ws1.Cells.Item(32, 6).Value = '黄礼闯'
$ws1.Cells.Item(32, 7).Value = 45281
$ws1.Cells.Item(32, 8).Value = '完成'
$ws1.Cells.Item(32, 9).Value = 'X药方-草药-单体-靶点'

# --- Apply date format to the new date cells (A21:A32, G21:G32) so they get a proper date style ---
$ws1.Range("A21:A32").NumberFormat = "m/d/yyyy"
$ws1.Range("G21:G32").NumberFormat = "m/d/yyyy"

# --- Re-apply the same date format to the pre-existing date columns so they migrate to the shared style ---
$ws1.Range("A2:A20").NumberFormat = "m/d/yyyy"
$ws1.Range("G2:G20").NumberFormat = "m/d/yyyy"
$ws2.Range("A2:A4").NumberFormat = "m/d/yyyy"
$ws2.Range("G2:G4").NumberFormat = "m/d/yyyy"
